# Commit: feat: add 2022-Q1 data
#
# The "2022-Q1" quarterly holdings sheet is new, and the running "总计"
# summary sheet gains a corresponding row. We reproduce the natural Excel
# workflow for this: duplicate the existing "总计" sheet (the duplicate is
# placed immediately after the original, exactly like Excel's "Move or
# Copy -> Create a copy"), turn the original into the new "2022-Q1" detail
# sheet, and turn the duplicate back into the refreshed "总计" sheet.

$wb = $excel.ActiveWorkbook

# --- duplicate "总计"; Excel places the copy right after the source ---
$src = $wb.Worksheets.Item("总计")
$src.Copy($null, $src)

$fundSheet  = $wb.Worksheets.Item("总计")
$totalSheet = $wb.Worksheets.Item("总计 (2)")

# free up the "总计" name on the original before renaming the copy to it
$fundSheet.Name = "2022-Q1"
$totalSheet.Name = "总计"

# =========================================================================
# "2022-Q1" sheet: per-fund holdings table (same layout as the other quarters)
# =========================================================================
$fundSheet.Range("A1:D6").ClearContents()

$fundSheet.Cells.Item(1,2).Value = "基金代码"
$fundSheet.Cells.Item(1,3).Value = "基金名称"
$fundSheet.Cells.Item(1,4).Value = "基金规模"
$fundSheet.Cells.Item(1,5).Value = "股票总仓位"
$fundSheet.Cells.Item(1,6).Value = "仓位占比"
$fundSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$fundSheet.Cells.Item(1,8).Value = "仓位排名"

# B1:D1 kept their header style from the original sheet; stamp the same
# style onto the newly-added E1:H1 header cells
$fundSheet.Cells.Item(1,3).Copy()
$fundSheet.Range("E1:H1").PasteSpecial(-4122)

$fundRows = @(
    @(0, "540008", "汇丰晋信低碳先锋股票", "96.71", "93.08", "8.90", "8.6072", 2),
    @(1, "003834", "华夏能源革新股票", "187.75", "93.26", "4.01", "7.5288", 9),
    @(2, "540003", "汇丰晋信动态策略混合A", "102.73", "91.87", "3.05", "3.1333", 9),
    @(3, "960003", "汇丰晋信动态策略混合H", "102.73", "91.87", "3.05", "3.1333", 9),
    @(4, "001643", "汇丰晋信智造先锋股票A", "29.09", "92.99", "8.83", "2.5686", 1),
    @(5, "011578", "汇丰晋信核心成长混合型证券投资基金A", "31.02", "91.97", "6.52", "2.0225", 2),
    @(6, "001644", "汇丰晋信智造先锋股票C", "10.91", "92.99", "8.83", "0.9634", 1),
    @(7, "003986", "申万菱信中证500指数优选增强A", "23.29", "92.19", "1.72", "0.4006", 7),
    @(8, "005477", "长安鑫禧灵活配置混合A", "5.70", "94.29", "5.23", "0.2981", 9),
    @(9, "011579", "汇丰晋信核心成长混合型证券投资基金C", "4.27", "91.97", "6.52", "0.2784", 2),
    @(10, "005343", "长安裕盛灵活配置混合A", "4.66", "94.21", "5.29", "0.2465", 6),
    @(11, "159870", "鹏华中证细分化工产业主题ETF", "8.45", "98.37", "2.70", "0.2282", 8),
    @(12, "562800", "嘉实中证稀有金属主题ETF", "6.36", "99.34", "3.57", "0.2271", 7),
    @(13, "005478", "长安鑫禧灵活配置混合C", "4.17", "94.29", "5.23", "0.2181", 9),
    @(14, "005344", "长安裕盛灵活配置混合C", "3.75", "94.21", "5.29", "0.1984", 6),
    @(15, "161039", "富国中证1000指数增强LOF", "21.72", "89.03", "0.77", "0.1672", 4),
    @(16, "512100", "南方中证1000ETF", "20.82", "99.21", "0.54", "0.1124", 3),
    @(17, "516020", "华宝中证细分化工产业主题ETF", "3.81", "98.92", "2.72", "0.1036", 8),
    @(18, "005313", "万家中证1000指数增强A", "9.01", "93.72", "1.10", "0.0991", 6),
    @(19, "159608", "广发中证稀有金属ETF", "2.39", "98.72", "3.51", "0.0839", 7),
    @(20, "516120", "富国中证细分化工产业主题ETF", "2.32", "99.12", "2.72", "0.0631", 8),
    @(21, "007794", "申万菱信中证500指数优选增强C", "3.53", "92.19", "1.72", "0.0607", 7),
    @(22, "014135", "中欧金安量化混合A", "9.43", "67.44", "0.63", "0.0594", 5),
    @(23, "005314", "万家中证1000指数增强C", "4.95", "93.72", "1.10", "0.0544", 6),
    @(24, "516220", "国泰中证细分化工产业主题ETF", "1.84", "98.95", "2.91", "0.0535", 7),
    @(25, "000646", "华润元大量化优选混合A", "0.80", "67.15", "6.53", "0.0522", 3),
    @(26, "290014", "泰信现代服务业混合", "0.73", "81.14", "6.68", "0.0488", 6),
    @(27, "001421", "南方量化成长股票", "1.70", "92.11", "2.54", "0.0432", 1),
    @(28, "002210", "创金合信量化多因子股票A", "3.31", "88.74", "1.16", "0.0384", 2),
    @(29, "290008", "泰信发展主题混合", "0.68", "81.03", "4.32", "0.0294", 9),
    @(30, "007827", "华润元大量化优选混合C", "0.33", "67.15", "6.53", "0.0215", 3),
    @(31, "001261", "中融新机遇灵活配置混合", "0.34", "93.08", "5.85", "0.0199", 6),
    @(32, "003646", "创金合信中证1000指数增强A", "1.42", "90.96", "1.19", "0.0169", 2),
    @(33, "003647", "创金合信中证1000指数增强C", "0.95", "90.96", "1.19", "0.0113", 2),
    @(34, "159845", "华夏中证1000ETF", "2.02", "97.95", "0.51", "0.0103", 4),
    @(35, "003865", "创金合信量化多因子股票C", "0.79", "88.74", "1.16", "0.0092", 2),
    @(36, "014136", "中欧金安量化混合C", "1.28", "67.44", "0.63", "0.0081", 5),
    @(37, "006486", "广发中证1000指数A", "1.26", "92.11", "0.49", "0.0062", 3),
    @(38, "159918", "嘉实中创400ETF", "0.75", "99.13", "0.79", "0.0059", 6),
    @(39, "004359", "创金合信量化核心混合A", "0.21", "93.31", "2.05", "0.0043", 7),
    @(40, "516300", "华泰柏瑞中证1000ETF", "0.49", "97.59", "0.53", "0.0026", 4),
    @(41, "006487", "广发中证1000指数C", "0.45", "92.11", "0.49", "0.0022", 3),
    @(42, "162413", "华宝中证1000指数", "0.37", "93.69", "0.51", "0.0019", 4),
    @(43, "006157", "财通量化核心优选混合", "0.09", "92.85", "1.72", "0.0015", 3),
    @(44, "001607", "英大策略优选混合A", "0.06", "93.31", "2.05", "0.0012", 7)
)

# B-G are stored as text (codes need leading zeros, decimals need exact
# trailing-zero formatting like "8.90"), so prefix with an apostrophe the
# same way a user typing into Excel would force text entry; H is numeric
$apos = "'"
$r = 2
foreach ($row in $fundRows) {
    $fundSheet.Cells.Item($r, 1).Value = $row[0]
    $fundSheet.Cells.Item($r, 2).Value = $apos + $row[1]
    $fundSheet.Cells.Item($r, 3).Value = $row[2]
    $fundSheet.Cells.Item($r, 4).Value = $apos + $row[3]
    $fundSheet.Cells.Item($r, 5).Value = $apos + $row[4]
    $fundSheet.Cells.Item($r, 6).Value = $apos + $row[5]
    $fundSheet.Cells.Item($r, 7).Value = $apos + $row[6]
    $fundSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# column A (row index) keeps the bold/boxed style from the original A2;
# stamp it onto every index cell down through row 46
$fundSheet.Cells.Item(2,1).Copy()
$fundSheet.Range("A2:A46").PasteSpecial(-4122)

# =========================================================================
# "总计" sheet: prepend the 2022-Q1 summary row, shifting the rest down
# =========================================================================
$totalRows = @(
    @(0, "2022-Q1", 45, 31.24),
    @(1, "2021-Q4", 36, 23.17),
    @(2, "2021-Q3", 64, 44.78),
    @(3, "2021-Q2", 31, 20.3),
    @(4, "2021-Q1", 11, 14.51),
    @(5, "2020-Q4", 28, 20.55)
)

$r = 2
foreach ($row in $totalRows) {
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# column A keeps the boxed index style; extend it to the newly-added row 7
$totalSheet.Cells.Item(2,1).Copy()
$totalSheet.Range("A2:A7").PasteSpecial(-4122)

Write-Output "2022-Q1 sheet added; 总计 sheet refreshed"
